# Generate Report for Handback
#
# The localization-status report is regenerated after the cae07396-fd5e-4559-
# ae4b-37b13a9513e0 file is handed back "in sync with en-US":
#   - its Status is flipped from "Ready for handoff" to
#     "Handed back: in sync with en-US" everywhere it is reported
#     (Overview sheet as well as the per-locale zh-cn / de-de sheets);
#   - on the per-locale sheets, the "Latest Handback DateTime" is refreshed
#     to the new handback timestamp, and the stale "Error Detail" (about the
#     handback file being out of date) is cleared since the file is now
#     current.

$wb = $excel.ActiveWorkbook

$status = "Handed back: in sync with en-US"

# --- Overview sheet: row 3 is the cae07396... file, reported once per locale ---
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("E3").Value = $status   # zh-cn column
$ws.Range("F3").Value = $status   # de-de column

# --- zh-cn sheet: row 3 is the cae07396... file ---
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("C3").Value = $status                  # Status
$ws.Range("K3").Value = "2016-08-17 14:49:11"     # Latest Handback DateTime
$ws.Range("P3").Value = ""                        # Error Detail (now resolved)
$ws.Columns.Item(16).AutoFit() | Out-Null         # Error Detail column no longer needs to be wide

# --- de-de sheet: row 3 is the cae07396... file ---
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("C3").Value = $status                  # Status
$ws.Range("K3").Value = "2016-08-17 14:49:19"     # Latest Handback DateTime
$ws.Range("P3").Value = ""                        # Error Detail (now resolved)
$ws.Columns.Item(16).AutoFit() | Out-Null         # Error Detail column no longer needs to be wide
